# Fill in the previously-blank C11:I11 scenario values on each worksheet
# ("Test 1" and "Test 2") with the updated R program / scenario figures.

$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 4.239895769806438
    "D11" = -0.09999999999999964
    "E11" = 0.8502771961109232
    "F11" = -0.07799999999999996
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 2.2179592315064
}

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
